$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 2102.0908
$ws.Range("I9").Value = 1134.6
$ws.Range("J9").Value = 2908.3333
$ws.Range("K9").Value = 1134.6
$ws.Range("L9").Value = 2908.3333
$ws.Range("M9").Value = -965.5999999999999
$ws.Range("N9").Value = -3246.3333
$ws.Range("H17").Value = 1081.5
$ws.Range("J17").Value = 1081.5
$ws.Range("L17").Value = 3244.5
$ws.Range("N17").Value = -3580.5
$ws.Range("H18").Value = 582.4
$ws.Range("J18").Value = 580
$ws.Range("L18").Value = 580
$ws.Range("N18").Value = -1148
$ws.Range("H62").Value = 166668590
$ws.Range("I62").Value = 194446290
$ws.Range("K62").Value = 194446290
$ws.Range("M62").Value = -194445666
$ws.Range("H65").Value = 166668590
$ws.Range("I65").Value = 194446290
$ws.Range("K65").Value = 972231450
$ws.Range("M65").Value = -972228330
$ws.Range("H106").Value = 2553.7727
$ws.Range("I106").Value = 2541.2104
$ws.Range("K106").Value = 2541.2104
$ws.Range("M106").Value = -1910.2104
$ws.Range("H116").Value = 15560.4
$ws.Range("I116").Value = 15622.667
$ws.Range("K116").Value = 15622.667
$ws.Range("M116").Value = -12180.667
$ws.Range("H132").Value = 1441.5714
$ws.Range("I132").Value = 1119.3334
$ws.Range("J132").Value = 2529.125
$ws.Range("K132").Value = 3358.0002
$ws.Range("L132").Value = 7587.375
$ws.Range("M132").Value = -828.0001999999999
$ws.Range("N132").Value = -12647.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 522.1818
$ws.Range("I4").Value = 467.44446
$ws.Range("J4").Value = 768.5
$ws.Range("K4").Value = 467.44446
$ws.Range("L4").Value = 768.5
$ws.Range("M4").Value = -351.44446
$ws.Range("N4").Value = -1000.5
$ws.Range("H74").Value = 2921.125
$ws.Range("I74").Value = 856.6875
$ws.Range("K74").Value = 856.6875
$ws.Range("M74").Value = 17.3125
$ws.Range("H77").Value = 2921.125
$ws.Range("I77").Value = 856.6875
$ws.Range("K77").Value = 4283.4375
$ws.Range("M77").Value = 84.5625
$ws.Range("H132").Value = 8224.538
$ws.Range("I132").Value = 5090.8237
$ws.Range("J132").Value = 14143.777
$ws.Range("K132").Value = 15272.4711
$ws.Range("L132").Value = 42431.331
$ws.Range("M132").Value = -12742.4711
$ws.Range("N132").Value = -47491.331

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2917.25
$ws.Range("I20").Value = 2721.2856
$ws.Range("J20").Value = 3113.2144
$ws.Range("K20").Value = 2721.2856
$ws.Range("L20").Value = 3113.2144
$ws.Range("M20").Value = -2474.2856
$ws.Range("N20").Value = -3607.2144
$ws.Range("H86").Value = 4254.5
$ws.Range("I86").Value = 4006
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 4006
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -2883
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 4254.5
$ws.Range("I89").Value = 4006
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 20030
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -14414
$ws.Range("N89").Value = -36232
$ws.Range("H94").Value = 1204.6923
$ws.Range("I94").Value = 1041.1
$ws.Range("J94").Value = 1750
$ws.Range("K94").Value = 1041.1
$ws.Range("L94").Value = 1750
$ws.Range("M94").Value = -590.0999999999999
$ws.Range("N94").Value = -2652
$ws.Range("H134").Value = 5897.7856
$ws.Range("I134").Value = 3269.5
$ws.Range("K134").Value = 9808.5
$ws.Range("M134").Value = -7273.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 452.5
$ws.Range("I22").Value = 470
$ws.Range("K22").Value = 470
$ws.Range("M22").Value = -120
$ws.Range("H31").Value = 2329909
$ws.Range("I31").Value = 1270.5834
$ws.Range("J31").Value = 3231317.5
$ws.Range("K31").Value = 1270.5834
$ws.Range("L31").Value = 3231317.5
$ws.Range("M31").Value = -975.5834
$ws.Range("N31").Value = -3231907.5
$ws.Range("H34").Value = 2329909
$ws.Range("I34").Value = 1270.5834
$ws.Range("J34").Value = 3231317.5
$ws.Range("K34").Value = 1270.5834
$ws.Range("L34").Value = 3231317.5
$ws.Range("M34").Value = -1068.5834
$ws.Range("N34").Value = -3231721.5
$ws.Range("H58").Value = 4938.5835
$ws.Range("I58").Value = 3750.7273
$ws.Range("J58").Value = 5943.6924
$ws.Range("K58").Value = 3750.7273
$ws.Range("L58").Value = 5943.6924
$ws.Range("M58").Value = -3547.7273
$ws.Range("N58").Value = -6349.6924
$ws.Range("H105").Value = 4394.467
$ws.Range("I105").Value = 6390.778
$ws.Range("J105").Value = 1400
$ws.Range("K105").Value = 6390.778
$ws.Range("L105").Value = 1400
$ws.Range("M105").Value = -4643.778
$ws.Range("N105").Value = -4894
$ws.Range("H134").Value = 5847.839
$ws.Range("I134").Value = 5930
$ws.Range("J134").Value = 5566.143
$ws.Range("K134").Value = 17790
$ws.Range("L134").Value = 16698.429
$ws.Range("M134").Value = -15255
$ws.Range("N134").Value = -21768.429
$ws.Range("H136").Value = 4938.5835
$ws.Range("I136").Value = 3750.7273
$ws.Range("J136").Value = 5943.6924
$ws.Range("K136").Value = 11252.1819
$ws.Range("L136").Value = 17831.0772
$ws.Range("M136").Value = -8702.1819
$ws.Range("N136").Value = -22931.0772

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 159
$ws.Range("I5").Value = 152.78572
$ws.Range("J5").Value = 202.5
$ws.Range("K5").Value = 458.35716
$ws.Range("L5").Value = 607.5
$ws.Range("M5").Value = -346.35716
$ws.Range("N5").Value = -831.5
$ws.Range("H34").Value = 2275.25
$ws.Range("J34").Value = 3715
$ws.Range("L34").Value = 11145
$ws.Range("N34").Value = -11313
$ws.Range("H68").Value = 120508.7
$ws.Range("J68").Value = 4000.2727
$ws.Range("L68").Value = 12000.8181
$ws.Range("N68").Value = -13622.8181
$ws.Range("H71").Value = 120508.7
$ws.Range("J71").Value = 4000.2727
$ws.Range("L71").Value = 36002.4543
$ws.Range("N71").Value = -44114.4543
$ws.Range("H75").Value = 624.1111
$ws.Range("I75").Value = 625.25
$ws.Range("J75").Value = 615
$ws.Range("K75").Value = 1875.75
$ws.Range("L75").Value = 1845
$ws.Range("M75").Value = -877.75
$ws.Range("N75").Value = -3841
$ws.Range("H78").Value = 624.1111
$ws.Range("I78").Value = 625.25
$ws.Range("J78").Value = 615
$ws.Range("K78").Value = 5627.25
$ws.Range("L78").Value = 5535
$ws.Range("M78").Value = -635.25
$ws.Range("N78").Value = -15519
$ws.Range("H129").Value = 1438.3334
$ws.Range("I129").Value = 332.33334
$ws.Range("J129").Value = 2544.3333
$ws.Range("K129").Value = 997.0000200000001
$ws.Range("L129").Value = 7632.999899999999
$ws.Range("M129").Value = 4002.99998
$ws.Range("N129").Value = -17632.9999
$ws.Range("H135").Value = 159
$ws.Range("I135").Value = 152.78572
$ws.Range("J135").Value = 202.5
$ws.Range("K135").Value = 1375.07148
$ws.Range("L135").Value = 1822.5
$ws.Range("M135").Value = 1159.92852
$ws.Range("N135").Value = -6892.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8228.9
$ws.Range("I70").Value = 7970.7144
$ws.Range("K70").Value = 7970.7144
$ws.Range("M70").Value = -7700.7144
$ws.Range("H73").Value = 8228.9
$ws.Range("I73").Value = 7970.7144
$ws.Range("K73").Value = 7970.7144
$ws.Range("M73").Value = -7034.7144
$ws.Range("H102").Value = 3157.625
$ws.Range("I102").Value = 3023
$ws.Range("K102").Value = 3023
$ws.Range("M102").Value = -1401
$ws.Range("H132").Value = 7474.9355
$ws.Range("I132").Value = 6444.357
$ws.Range("J132").Value = 8323.647000000001
$ws.Range("K132").Value = 19333.071
$ws.Range("L132").Value = 24970.941
$ws.Range("M132").Value = -16803.071
$ws.Range("N132").Value = -30030.941

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1506.0476
$ws.Range("I16").Value = 1497.6316
$ws.Range("J16").Value = 1586
$ws.Range("K16").Value = 1497.6316
$ws.Range("L16").Value = 1586
$ws.Range("M16").Value = -1327.6316
$ws.Range("N16").Value = -1926
$ws.Range("H46").Value = 5419.6816
$ws.Range("J46").Value = 6856.0625
$ws.Range("L46").Value = 6856.0625
$ws.Range("N46").Value = -7232.0625
$ws.Range("H93").Value = 3286.5881
$ws.Range("I93").Value = 7478.6665
$ws.Range("K93").Value = 7478.6665
$ws.Range("M93").Value = -6230.6665
$ws.Range("H122").Value = 7666.3335
$ws.Range("I122").Value = 9499.5
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 28498.5
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -26048.5
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 4137.927
$ws.Range("I132").Value = 3672.158
$ws.Range("J132").Value = 4540.1816
$ws.Range("K132").Value = 11016.474
$ws.Range("L132").Value = 13620.5448
$ws.Range("M132").Value = -8486.474
$ws.Range("N132").Value = -18680.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 24946.8
$ws.Range("I52").Value = 24941
$ws.Range("K52").Value = 24941
$ws.Range("M52").Value = -24715
$ws.Range("H135").Value = 90714.5
$ws.Range("J135").Value = 90714.5
$ws.Range("L135").Value = 90714.5
$ws.Range("N135").Value = -100854.5
